$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(13, 14, 17, 18, 19, 20, 21, 22)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 9).Value = "Rejected"
    $ws.Cells.Item($r, 10).Value = "Nil"
}

$ws.Range("J15").Select()
